# Updated cryptos list on Thu Dec 28 20:11:01 UTC 2023 with GitHub Actions
#
# Helper: assign a value that must remain plain text (not get auto-coerced
# into a number by Excel) while keeping the cell's original (unset / default)
# style, matching the source workbook's inlineStr cells.
function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "42.668.74"
Set-TextValue $ws "E2" "  -1.38%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "2.357.69"
Set-TextValue $ws "E3" "  +0.26%  "

# Row 4 - TetherUSD
Set-TextValue $ws "E4" "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "328.52"
Set-TextValue $ws "E5" "  +5.05%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "100.84"
Set-TextValue $ws "E6" "  -7.95%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.636"
Set-TextValue $ws "E7" "  -1.12%  "

# Row 8 - USDC
Set-TextValue $ws "E8" "  +0.13%  "

# Row 9 - Cardano
Set-TextValue $ws "D9" "0.625"
Set-TextValue $ws "E9" "  -2.12%  "

# Row 10 - Avalanche
Set-TextValue $ws "D10" "39.76"
Set-TextValue $ws "E10" "  -7.71%  "

# Row 11 - Dogecoin
Set-TextValue $ws "D11" "0.0921"
Set-TextValue $ws "E11" "  -1.77%  "

# Row 12 - Polkadot
Set-TextValue $ws "D12" "8.43"
Set-TextValue $ws "E12" "  -4.68%  "

# Row 13 - Polygon
Set-TextValue $ws "D13" "1.00"
Set-TextValue $ws "E13" "  -4.08%  "

# Row 14 - TRON
Set-TextValue $ws "E14" "  +0.21%  "

# Row 15 - Chainlink
Set-TextValue $ws "D15" "16.46"
Set-TextValue $ws "E15" "  +0.24%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D16" "2.714.62"
Set-TextValue $ws "E16" "  +0.45%  "

# Row 17 - WrappedEther
Set-TextValue $ws "D17" "2.357.47"
Set-TextValue $ws "E17" "  -4.86%  "

# Row 18 - Uniswap
Set-TextValue $ws "D18" "8.11"
Set-TextValue $ws "E18" "  +11.99%  "

# Row 19 - WrappedBTC
Set-TextValue $ws "D19" "42.666.53"
Set-TextValue $ws "E19" "  -1.36%  "

# Row 20 - ShibaInu
Set-TextValue $ws "E20" "  -1.94%  "

# Row 21 - Litecoin
Set-TextValue $ws "D21" "76.21"
Set-TextValue $ws "E21" "  +1.02%  "

# Row 22 - PancakeSwap
Set-TextValue $ws "D22" "3.73"
Set-TextValue $ws "E22" "  +8.30%  "

# Row 23 - BitcoinCash
Set-TextValue $ws "D23" "269.88"
Set-TextValue $ws "E23" "  +5.85%  "

# Row 24 & 25 swap: InternetComputer(DFINITY) <-> ImmutableX
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D24" "2.31"
Set-TextValue $ws "E24" "  -10.22%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D25" "10.14"
Set-TextValue $ws "E25" "  +11.16%  "

# Row 26 - Dai
Set-TextValue $ws "E26" "  -0.10%  "

# Row 27 - Cosmos
Set-TextValue $ws "D27" "11.49"
Set-TextValue $ws "E27" "  -4.60%  "

# Row 28 - EthereumClassic
Set-TextValue $ws "D28" "23.07"
Set-TextValue $ws "E28" "  +3.19%  "

# Row 29 - Toncoin
Set-TextValue $ws "D29" "2.20"
Set-TextValue $ws "E29" "  -2.11%  "

# Row 30 - Monero
Set-TextValue $ws "D30" "176.44"
Set-TextValue $ws "E30" "  +1.43%  "

# Row 31 - WEMIXToken
Set-TextValue $ws "E31" "  -2.38%  "

# Row 32 - Hedera
Set-TextValue $ws "D32" "0.0901"
Set-TextValue $ws "E32" "  -3.05%  "

# Row 33 - InjectiveProtocol
Set-TextValue $ws "D33" "35.45"
Set-TextValue $ws "E33" "  -9.60%  "

# Row 34 - Filecoin
Set-TextValue $ws "D34" "6.07"
Set-TextValue $ws "E34" "  +0.46%  "

# Row 35 - Stellar
Set-TextValue $ws "D35" "0.132"
Set-TextValue $ws "E35" "  -0.13%  "

# Row 36 - RenderToken
Set-TextValue $ws "D36" "4.61"
Set-TextValue $ws "E36" "  -7.46%  "

# Row 37 & 38 swap: LidoDAOToken <-> VeChain
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D37" "0.0359"
Set-TextValue $ws "E37" "  -4.64%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D38" "2.94"
Set-TextValue $ws "E38" "  +8.96%  "

# Row 39 - Kaspa
Set-TextValue $ws "E39" "  +1.08%  "

# Row 40 - NEARProtocol
Set-TextValue $ws "D40" "3.79"
Set-TextValue $ws "E40" "  -8.37%  "

# Row 41 - ARBITRUM
Set-TextValue $ws "D41" "1.51"
Set-TextValue $ws "E41" "  +1.86%  "

# Row 42 - Algorand
Set-TextValue $ws "D42" "0.236"
Set-TextValue $ws "E42" "  +1.07%  "

# Row 43 - MultiversX
Set-TextValue $ws "D43" "69.91"
Set-TextValue $ws "E43" "  -4.03%  "

# Row 44 - FirstDigitalUSD
Set-TextValue $ws "E44" "  -0.01%  "

# Row 45 & 46 swap: Aave <-> BitcoinSV
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws "D45" "92.58"
Set-TextValue $ws "E45" "  +32.29%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D46" "118.80"
Set-TextValue $ws "E46" "  +6.91%  "

# Row 47 - Celestia
Set-TextValue $ws "D47" "11.77"
Set-TextValue $ws "E47" "  -7.72%  "

# Row 48 - THORChain
Set-TextValue $ws "D48" "5.50"
Set-TextValue $ws "E48" "  -2.03%  "

# Row 49 - FraxShare
Set-TextValue $ws "D49" "9.19"
Set-TextValue $ws "E49" "  -1.06%  "

# Row 50 - TrustWalletToken
Set-TextValue $ws "E50" "  -3.05%  "

# Row 51 - Maker
Set-TextValue $ws "D51" "1.569.94"
Set-TextValue $ws "E51" "  +5.06%  "
